$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q3" right before the current "2022-Q2"
#    sheet. Copying the existing "2022-Q2" sheet gives us the right header
#    row / column styles "for free", then we overwrite the data rows with
#    the new quarter's figures and drop the extra (now-stale) rows.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Columns B, D, E, F, G hold numeric-looking values that are stored as text
# in the source data (e.g. fund codes with leading zeros, percentages kept
# as plain strings) - force text format so Excel doesn't silently convert
# them to numbers when we assign the string values.
$wsQ3.Range("B2:B4").NumberFormat = "@"
$wsQ3.Range("D2:G4").NumberFormat = "@"

$wsQ3.Range("B2").Value = "006022"
$wsQ3.Range("C2").Value = "富国大盘价值量化精选混合A"
$wsQ3.Range("D2").Value = "5.68"
$wsQ3.Range("E2").Value = "91.40"
$wsQ3.Range("F2").Value = "1.89"
$wsQ3.Range("G2").Value = "0.1074"
$wsQ3.Range("H2").Value = 8

$wsQ3.Range("B3").Value = "001068"
$wsQ3.Range("C3").Value = "华融新锐灵活配置混合"
$wsQ3.Range("D3").Value = "0.21"
$wsQ3.Range("E3").Value = "53.47"
$wsQ3.Range("F3").Value = "2.70"
$wsQ3.Range("G3").Value = "0.0057"
$wsQ3.Range("H3").Value = 6

$wsQ3.Range("B4").Value = "014181"
$wsQ3.Range("C4").Value = "富国大盘价值量化精选混合C"
$wsQ3.Range("D4").Value = "0.11"
$wsQ3.Range("E4").Value = "91.40"
$wsQ3.Range("F4").Value = "1.89"
$wsQ3.Range("G4").Value = "0.0021"
$wsQ3.Range("H4").Value = 8

# The text values are now locked in as strings - drop the temporary "@"
# number format again so the cells fall back to the plain/default style
# (matching the unstyled data cells used everywhere else in this workbook).
$wsQ3.Range("B2:B4").Style = "Normal"
$wsQ3.Range("D2:G4").Style = "Normal"

# The copied sheet still has the old rows 5,6,7 - delete them, the new
# quarter only has 3 funds.
$wsQ3.Rows.Item(7).Delete()
$wsQ3.Rows.Item(6).Delete()
$wsQ3.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    after the header, push the existing quarters down one row, and
#    renumber the running index in column A.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Copy the formatting (border/alignment style) of the last existing data
# row down into the new row 6.
$wsTotal.Range("A5").Copy($wsTotal.Range("A6"))

$wsTotal.Range("B6").Value = $wsTotal.Range("B5").Value2
$wsTotal.Range("C6").Value = $wsTotal.Range("C5").Value2
$wsTotal.Range("D6").Value = $wsTotal.Range("D5").Value2

$wsTotal.Range("B5").Value = $wsTotal.Range("B4").Value2
$wsTotal.Range("C5").Value = $wsTotal.Range("C4").Value2
$wsTotal.Range("D5").Value = $wsTotal.Range("D4").Value2

$wsTotal.Range("B4").Value = $wsTotal.Range("B3").Value2
$wsTotal.Range("C4").Value = $wsTotal.Range("C3").Value2
$wsTotal.Range("D4").Value = $wsTotal.Range("D3").Value2

$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value2
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value2
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value2

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.12

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
